# Append the 2023 row to the "data" sheet (annual stats for the 2023 report).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Column A holds the year as text (matching the existing "year" column,
# which is stored as shared-string text, not numbers). Enter it through a
# text formula and convert the result to a static value so it lands in the
# workbook the same way the other year cells did.
$ws.Range("A23").Formula = "=""2023"""
$ws.Range("A23").Copy() | Out-Null
$ws.Range("A23").PasteSpecial(-4163) | Out-Null  # xlPasteValues
$excel.CutCopyMode = 0

$ws.Range("B23").Value = 246344
$ws.Range("C23").Value = 7.9
$ws.Range("D23").Value = 4807
$ws.Range("E23").Value = 38081

# Reflect where the user's selection ended up after adding the new row.
$ws.Range("E22").Select() | Out-Null
